$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. "1.000", "27.074.49")
# keep their exact literal representation instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.074.49"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "1.825.66"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -1.29%  "
$ws.Range("D5").Value = "311.42"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "0.4228"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "0.3679"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "0.07232"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "0.8439"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").Value = "20.77"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").Value = "1.813.08"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "6.675"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D15").Value = "5.297"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").Value = "89.84"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "0.000008750"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "27.102.98"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "5.149"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").Value = "2.049.43"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "1.985"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "151.72"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "2.248"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").Value = "18.30"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "5.265"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "116.95"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "0.08738"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "1.178"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "0.7371"
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("D34").Value = "4.423"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "2.881"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "1.093"
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("D38").Value = "0.05260"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").Value = "7.319"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").Value = "2.871"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "0.1687"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "0.5066"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "8.572"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "10.54"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "106.08"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "0.4725"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "1.942"
$ws.Range("E48").Value = "  +5.54%  "
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").Value = "0.06326"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").Value = "1.653"
$ws.Range("E51").Value = "  -2.23%  "
